$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.434.17"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "1.691.42"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5531"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.009"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2715"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07610"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "1.689.54"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5837"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008467"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "26.530.85"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.960"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.259"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "150.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("E25").Value = "  +9.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.951"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.416"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.66%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06309"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.593"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.593"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.678"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.050"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6250"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.406"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.252"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.124.21"
$ws.Range("E39").Value = "  +2.83%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01644"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "1.841.78"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.220"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05283"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4301"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
